$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4444.4443
$ws.Cells.Item(17, 10).Value = 5599.8
$ws.Cells.Item(17, 12).Value = 16799.4
$ws.Cells.Item(17, 14).Value = -17135.4
$ws.Cells.Item(62, 8).Value = 4195.25
$ws.Cells.Item(62, 9).Value = 3380.875
$ws.Cells.Item(62, 11).Value = 3380.875
$ws.Cells.Item(62, 13).Value = -2756.875
$ws.Cells.Item(64, 8).Value = 6833.3335
$ws.Cells.Item(64, 9).Value = 6833.3335
$ws.Cells.Item(64, 11).Value = 6833.3335
$ws.Cells.Item(64, 13).Value = -6585.3335
$ws.Cells.Item(65, 8).Value = 4195.25
$ws.Cells.Item(65, 9).Value = 3380.875
$ws.Cells.Item(65, 11).Value = 16904.375
$ws.Cells.Item(65, 13).Value = -13784.375
$ws.Cells.Item(67, 8).Value = 6833.3335
$ws.Cells.Item(67, 9).Value = 6833.3335
$ws.Cells.Item(67, 11).Value = 6833.3335
$ws.Cells.Item(67, 13).Value = -5975.3335
$ws.Cells.Item(116, 8).Value = 4585
$ws.Cells.Item(116, 9).Value = 1999
$ws.Cells.Item(116, 11).Value = 1999
$ws.Cells.Item(116, 13).Value = 1443
$ws.Cells.Item(127, 8).Value = 2013.1
$ws.Cells.Item(127, 9).Value = 1459
$ws.Cells.Item(127, 11).Value = 4377
$ws.Cells.Item(127, 13).Value = 583
$ws.Cells.Item(129, 8).Value = 2146.7778
$ws.Cells.Item(129, 9).Value = 806.3333
$ws.Cells.Item(129, 10).Value = 3822.3333
$ws.Cells.Item(129, 11).Value = 2418.9999
$ws.Cells.Item(129, 12).Value = 11466.9999
$ws.Cells.Item(129, 13).Value = 2581.0001
$ws.Cells.Item(129, 14).Value = -21466.9999
$ws.Cells.Item(132, 8).Value = 3628.5715
$ws.Cells.Item(132, 9).Value = 4066.8333
$ws.Cells.Item(132, 11).Value = 12200.4999
$ws.Cells.Item(132, 13).Value = -9670.499899999999
$ws.Cells.Item(133, 8).Value = 38999
$ws.Cells.Item(133, 10).Value = 38999
$ws.Cells.Item(133, 12).Value = 38999
$ws.Cells.Item(133, 14).Value = -49119
$ws.Cells.Item(135, 8).Value = 2311.7778
$ws.Cells.Item(135, 9).Value = 1115.2858
$ws.Cells.Item(135, 10).Value = 6499.5
$ws.Cells.Item(135, 11).Value = 10037.5722
$ws.Cells.Item(135, 12).Value = 58495.5
$ws.Cells.Item(135, 13).Value = -7502.572200000001
$ws.Cells.Item(135, 14).Value = -63565.5
$ws.Cells.Item(138, 8).Value = 1400.409
$ws.Cells.Item(138, 9).Value = 535.44446
$ws.Cells.Item(138, 11).Value = 1606.33338
$ws.Cells.Item(138, 13).Value = 3533.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(33, 8).Value = 4500.5
$ws.Cells.Item(33, 9).Value = 4500.5
$ws.Cells.Item(33, 11).Value = 4500.5
$ws.Cells.Item(33, 13).Value = -4171.5
$ws.Cells.Item(97, 8).Value = 229.14285
$ws.Cells.Item(97, 9).Value = 229.14285
$ws.Cells.Item(97, 11).Value = 229.14285
$ws.Cells.Item(97, 13).Value = 266.85715

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1252.3478
$ws.Cells.Item(86, 9).Value = 1384.75
$ws.Cells.Item(86, 10).Value = 949.7143
$ws.Cells.Item(86, 11).Value = 1384.75
$ws.Cells.Item(86, 12).Value = 949.7143
$ws.Cells.Item(86, 13).Value = -261.75
$ws.Cells.Item(86, 14).Value = -3195.7143
$ws.Cells.Item(89, 8).Value = 1252.3478
$ws.Cells.Item(89, 9).Value = 1384.75
$ws.Cells.Item(89, 10).Value = 949.7143
$ws.Cells.Item(89, 11).Value = 6923.75
$ws.Cells.Item(89, 12).Value = 4748.5715
$ws.Cells.Item(89, 13).Value = -1307.75
$ws.Cells.Item(89, 14).Value = -15980.5715
$ws.Cells.Item(94, 8).Value = 2029.1333
$ws.Cells.Item(94, 9).Value = 2029.1333
$ws.Cells.Item(94, 11).Value = 2029.1333
$ws.Cells.Item(94, 13).Value = -1578.1333
$ws.Cells.Item(99, 8).Value = 1174.625
$ws.Cells.Item(99, 9).Value = 1128.1428
$ws.Cells.Item(99, 11).Value = 1128.1428
$ws.Cells.Item(99, 13).Value = 369.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1500
$ws.Cells.Item(34, 8).Value = 1500
$ws.Cells.Item(41, 8).Value = 23750
$ws.Cells.Item(41, 10).Value = 23750
$ws.Cells.Item(41, 12).Value = 23750
$ws.Cells.Item(41, 14).Value = -24606
$ws.Cells.Item(50, 8).Value = 29998
$ws.Cells.Item(50, 10).Value = 29998
$ws.Cells.Item(50, 12).Value = 29998
$ws.Cells.Item(50, 14).Value = -31248
$ws.Cells.Item(59, 8).Value = 32498.25
$ws.Cells.Item(59, 9).Value = 20000
$ws.Cells.Item(59, 10).Value = 34997.9
$ws.Cells.Item(59, 11).Value = 20000
$ws.Cells.Item(59, 12).Value = 34997.9
$ws.Cells.Item(59, 13).Value = -18855
$ws.Cells.Item(59, 14).Value = -37287.9
$ws.Cells.Item(60, 8).Value = 22998
$ws.Cells.Item(60, 10).Value = 24997.5
$ws.Cells.Item(60, 12).Value = 24997.5
$ws.Cells.Item(60, 14).Value = -26019.5
$ws.Cells.Item(68, 8).Value = 36562.5
$ws.Cells.Item(68, 9).Value = 12500
$ws.Cells.Item(68, 11).Value = 12500
$ws.Cells.Item(68, 13).Value = -11751
$ws.Cells.Item(71, 8).Value = 36562.5
$ws.Cells.Item(71, 9).Value = 12500
$ws.Cells.Item(71, 11).Value = 37500
$ws.Cells.Item(71, 13).Value = -33756
$ws.Cells.Item(141, 8).Value = 139999.83
$ws.Cells.Item(141, 9).Value = 90000
$ws.Cells.Item(141, 11).Value = 90000
$ws.Cells.Item(141, 13).Value = -84820

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 159.85715
$ws.Cells.Item(38, 9).Value = 209
$ws.Cells.Item(38, 11).Value = 627
$ws.Cells.Item(38, 13).Value = -280
$ws.Cells.Item(132, 8).Value = 2995
$ws.Cells.Item(132, 9).Value = 2995
$ws.Cells.Item(132, 11).Value = 26955
$ws.Cells.Item(132, 13).Value = -24425
$ws.Cells.Item(140, 8).Value = 1967.9
$ws.Cells.Item(140, 9).Value = 1967.9
$ws.Cells.Item(140, 11).Value = 5903.700000000001
$ws.Cells.Item(140, 13).Value = -723.7000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1321.25
$ws.Cells.Item(97, 9).Value = 1124.3334
$ws.Cells.Item(97, 10).Value = 1439.4
$ws.Cells.Item(97, 11).Value = 1124.3334
$ws.Cells.Item(97, 12).Value = 1439.4
$ws.Cells.Item(97, 13).Value = -628.3334
$ws.Cells.Item(97, 14).Value = -2431.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 9395.556
$ws.Cells.Item(40, 9).Value = 8639.75
$ws.Cells.Item(40, 11).Value = 8639.75
$ws.Cells.Item(40, 13).Value = -8503.75
$ws.Cells.Item(46, 8).Value = 2731.5264
$ws.Cells.Item(46, 10).Value = 4139.8
$ws.Cells.Item(46, 12).Value = 4139.8
$ws.Cells.Item(46, 14).Value = -4515.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 3500
$ws.Cells.Item(21, 9).Value = 3500
$ws.Cells.Item(21, 11).Value = 3500
$ws.Cells.Item(21, 13).Value = -3265
$ws.Cells.Item(30, 8).Value = 2595
$ws.Cells.Item(30, 10).Value = 2595
$ws.Cells.Item(30, 12).Value = 2595
$ws.Cells.Item(30, 14).Value = -2809
$ws.Cells.Item(31, 8).Value = 29500
$ws.Cells.Item(31, 9).Value = 20000
$ws.Cells.Item(31, 11).Value = 20000
$ws.Cells.Item(31, 13).Value = -19652
$ws.Cells.Item(35, 8).Value = 3500
$ws.Cells.Item(35, 9).Value = 3500
$ws.Cells.Item(35, 11).Value = 3500
$ws.Cells.Item(35, 13).Value = -3210
$ws.Cells.Item(96, 8).Value = 1598.75
$ws.Cells.Item(96, 9).Value = 1465
$ws.Cells.Item(96, 10).Value = 2000
$ws.Cells.Item(96, 11).Value = 1465
$ws.Cells.Item(96, 12).Value = 2000
$ws.Cells.Item(96, 13).Value = -92
$ws.Cells.Item(96, 14).Value = -4746
$ws.Cells.Item(101, 8).Value = 51000
$ws.Cells.Item(101, 10).Value = 51000
$ws.Cells.Item(101, 12).Value = 51000
$ws.Cells.Item(101, 14).Value = -57490
$ws.Cells.Item(135, 8).Value = 44249.5
$ws.Cells.Item(135, 10).Value = 44249.5
$ws.Cells.Item(135, 12).Value = 44249.5
$ws.Cells.Item(135, 14).Value = -54389.5
